$wb = $excel.ActiveWorkbook

# --- Sheet "Batter" ---
$ws = $wb.Worksheets.Item("Batter")

$ws.Range("C2").Value = 4.443236247967362
$ws.Range("E2").Value = -17.06722924051244

$ws.Range("C3").Value = 3.927950448455748
$ws.Range("E3").Value = -13.11967645362222

$ws.Range("C4").Value = 4.177346225815146
$ws.Range("E4").Value = -14.96958711751675

$ws.Range("C5").Value = 4.00401819782806
$ws.Range("D5").Value = 0.9884672369025292
$ws.Range("E5").Value = -13.67184835148691

$ws.Range("C6").Value = 4.260782067575361
$ws.Range("E6").Value = -15.6138921449346

$ws.Range("C7").Value = 3.924361314618307
$ws.Range("D7").Value = 0.9175794647558174
$ws.Range("E7").Value = -13.09388475581119

$ws.Range("C8").Value = 3.98931807735444
$ws.Range("D8").Value = 0.9846448435433534
$ws.Range("E8").Value = -13.5643153613973

# --- Sheet "Pitcher" ---
$ws2 = $wb.Worksheets.Item("Pitcher")

$ws2.Range("C2").Value = 32.39192243400829
$ws2.Range("E2").Value = -327.4602300122007

$ws2.Range("C3").Value = 9.352564995263009
$ws2.Range("E3").Value = -26.38235616466474

$ws2.Range("C4").Value = 7.807376003184183
$ws2.Range("E4").Value = -18.08180863123216

$ws2.Range("C5").Value = 11.47327868378975
$ws2.Range("D5").Value = 0.9913980763381751
$ws2.Range("E5").Value = -40.20827455014338

$ws2.Range("C6").Value = 22.59497503356168
$ws2.Range("E6").Value = -158.8207177224199

$ws2.Range("C7").Value = 10.64354763602054
$ws2.Range("D7").Value = 0.913653446242111
$ws2.Range("E7").Value = -34.46354624278319

$ws2.Range("C8").Value = 11.523215889036
$ws2.Range("D8").Value = 0.9878060360374233
$ws2.Range("E8").Value = -40.56777150708511
